$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '38.134.91'
$ws.Range('E2').Value = '  +2.90%  '
Set-TextValue $ws 'D3' '2.053.98'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('E4').Value = '  +0.17%  '
Set-TextValue $ws 'D5' '230.16'
$ws.Range('E5').Value = '  +1.16%  '
Set-TextValue $ws 'D6' '0.617'
$ws.Range('E6').Value = '  +1.14%  '
Set-TextValue $ws 'D7' '60.92'
$ws.Range('E7').Value = '  +9.19%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +2.62%  '
Set-TextValue $ws 'D10' '0.0804'
$ws.Range('E10').Value = '  +3.31%  '
$ws.Range('E11').Value = '  +2.12%  '
Set-TextValue $ws 'D12' '14.81'
$ws.Range('E12').Value = '  +4.13%  '
Set-TextValue $ws 'D13' '2.363.78'
$ws.Range('E13').Value = '  +2.05%  '
Set-TextValue $ws 'D14' '21.06'
$ws.Range('E14').Value = '  +5.67%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D15' '5.34'
$ws.Range('E15').Value = '  +3.38%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D16' '0.758'
$ws.Range('E16').Value = '  +2.94%  '
Set-TextValue $ws 'D17' '2.055.63'
$ws.Range('E17').Value = '  +1.74%  '
Set-TextValue $ws 'D18' '38.156.00'
$ws.Range('E18').Value = '  +2.97%  '
Set-TextValue $ws 'D19' '6.28'
$ws.Range('E19').Value = '  +2.54%  '
Set-TextValue $ws 'D20' '69.77'
$ws.Range('E20').Value = '  +1.28%  '
Set-TextValue $ws 'D21' '0.0₃0832'
Set-TextValue $ws 'D22' '226.01'
$ws.Range('E22').Value = '  +1.40%  '
Set-TextValue $ws 'D23' '0.999'
$ws.Range('E23').Value = '  +0.01%  '
Set-TextValue $ws 'D24' '2.43'
$ws.Range('E24').Value = '  +0.02%  '
Set-TextValue $ws 'D25' '2.22'
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('E26').Value = '  +3.06%  '
Set-TextValue $ws 'D27' '165.52'
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E28').Value = '  +4.82%  '
Set-TextValue $ws 'D29' '19.04'
$ws.Range('E29').Value = '  +2.20%  '
Set-TextValue $ws 'D30' '1.31'
$ws.Range('E30').Value = '  +0.47%  '
Set-TextValue $ws 'D31' '0.120'
$ws.Range('E31').Value = '  +2.27%  '
$ws.Range('E32').Value = '  +2.60%  '
Set-TextValue $ws 'D33' '4.57'
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D34' '2.05'
$ws.Range('E34').Value = '  +9.20%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D35' '0.0604'
$ws.Range('E35').Value = '  +0.47%  '
Set-TextValue $ws 'D36' '6.30'
$ws.Range('E36').Value = '  +15.02%  '
Set-TextValue $ws 'D37' '2.30'
$ws.Range('E37').Value = '  -1.85%  '
Set-TextValue $ws 'D38' '3.28'
$ws.Range('E38').Value = '  +3.89%  '
$ws.Range('E39').Value = '  +0.25%  '
Set-TextValue $ws 'D40' '1.516.33'
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D41' '0.0217'
$ws.Range('E41').Value = '  +2.03%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D42' '97.41'
$ws.Range('E42').Value = '  +3.38%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D43' '16.89'
$ws.Range('E43').Value = '  +4.36%  '
$ws.Range('E44').Value = '  +3.58%  '
Set-TextValue $ws 'D45' '0.0927'
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('E46').Value = '  +1.85%  '
Set-TextValue $ws 'D47' '4.03'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('E48').Value = '  +1.47%  '
Set-TextValue $ws 'D49' '2.97'
$ws.Range('E49').Value = '  +1.68%  '
Set-TextValue $ws 'D50' '7.05'
$ws.Range('E50').Value = '  +0.15%  '
Set-TextValue $ws 'D51' '2.253.76'
$ws.Range('E51').Value = '  +2.22%  '
